$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Awesome web design"
$ws.Range("A3").Value = "Terrible web applications."
$ws.Range("A4").Value = "Great and accuracy models."

$ws.Range("A1").Select()
